$x = [char]0x00D7
$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# (Row, Col, oldA, oldB, oldC, newA, newB, newC)
$edits = @(
    ,@(1, 1, "58", "59", "3422", "89", "66", "5874")
    ,@(1, 2, "44", "65", "2860", "71", "80", "5680")
    ,@(1, 3, "50", "53", "2650", "82", "33", "2706")
    ,@(1, 4, "35", "65", "2275", "23", "76", "1748")
    ,@(1, 5, "41", "32", "1312", "41", "45", "1845")
    ,@(5, 1, "77", "34", "2618", "70", "41", "2870")
    ,@(5, 2, "17", "41", "697", "15", "11", "165")
    ,@(5, 3, "99", "76", "7524", "98", "52", "5096")
    ,@(5, 4, "90", "76", "6840", "33", "65", "2145")
    ,@(5, 5, "67", "94", "6298", "73", "26", "1898")
    ,@(10, 1, "56", "54", "3024", "22", "12", "264")
    ,@(10, 2, "50", "86", "4300", "61", "37", "2257")
    ,@(10, 3, "11", "59", "649", "13", "66", "858")
    ,@(10, 4, "59", "50", "2950", "76", "90", "6840")
    ,@(10, 5, "30", "60", "1800", "65", "68", "4420")
    ,@(15, 1, "42", "53", "2226", "97", "59", "5723")
    ,@(15, 2, "70", "23", "1610", "89", "73", "6497")
    ,@(15, 3, "51", "50", "2550", "70", "92", "6440")
    ,@(15, 4, "41", "34", "1394", "42", "72", "3024")
    ,@(15, 5, "87", "89", "7743", "16", "95", "1520")
    ,@(20, 1, "50", "79", "3950", "97", "50", "4850")
    ,@(20, 2, "22", "85", "1870", "87", "67", "5829")
    ,@(20, 3, "28", "30", "840", "74", "52", "3848")
    ,@(20, 4, "21", "15", "315", "50", "53", "2650")
    ,@(20, 5, "95", "91", "8645", "53", "88", "4664")
)

foreach ($e in $edits) {
    $row = $e[0]
    $col = $e[1]
    $old = "{0}{1}{2}={3}" -f $e[2], $x, $e[3], $e[4]
    $new = "{0}{1}{2}={3}" -f $e[5], $x, $e[6], $e[7]
    $cell = $t.Cell($row, $col)
    $rng = $cell.Range
    $ok = $rng.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
    if (-not $ok) {
        Write-Host "FAILED row=$row col=$col old=$old"
    }
}

# Title line update
$result = $d.Content.Find.Execute("2024-10-19 Saturday", $true, $false, $false, $false, $false, $true, 1, $false, "2024-10-20 Sunday", 2)
if (-not $result) {
    Write-Host "FAILED title replace"
}

Write-Host "Done"
